$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.951.68"
$ws.Range("D3").Value = "2.219.17"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'292.22"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "'87.00"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("D9").Value = "'0.467"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").Value = "'30.49"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0780"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "'50.13"
$ws.Range("E12").Value = "  +5.17%  "
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "2.563.38"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.302.11"
$ws.Range("E16").Value = "  +4.03%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'13.77"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "39.869.83"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "'11.18"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "'65.53"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "'236.87"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "'1.82"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("E28").Value = "  +7.56%  "
$ws.Range("D29").Value = "'23.24"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").Value = "'157.59"
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D36").Value = "'0.0713"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D39").Value = "'0.0986"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").Value = "'15.17"
$ws.Range("E41").Value = "  -7.22%  "
$ws.Range("D42").Value = "2.090.01"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").Value = "'3.72"
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "'17.90"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("D46").Value = "'9.80"
$ws.Range("E46").Value = "  -3.30%  "
$ws.Range("D47").Value = "'2.01"
$ws.Range("E47").Value = "  -8.13%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").Value = "2.434.77"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("E51").Value = "  +2.17%  "
